$d = $word.ActiveDocument

# Locate the existing "{{ tipo_acao }}" placeholder and insert a new
# "{{ num_acao }}" placeholder (prefixed with " nº ") right after it,
# before the comma that starts "...bem como para, agindo em comum...".
$rng = $d.Content
$found = $rng.Find.Execute("tipo_acao }}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertPos = $rng.End
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertAfter(" nº {{ num_acao }}")
}
